$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-22 Wednesday" "2025-10-23 Thursday"

Replace-Text "559÷8=69, 7" "401÷6=66, 5"
Replace-Text "663÷6=110, 3" "913÷6=152, 1"
Replace-Text "288÷9=32, 0" "409÷6=68, 1"
Replace-Text "682÷5=136, 2" "198÷4=49, 2"
Replace-Text "639÷9=71, 0" "939÷3=313, 0"

Replace-Text "673÷7=96, 1" "548÷5=109, 3"
Replace-Text "534÷9=59, 3" "274÷4=68, 2"
Replace-Text "359÷8=44, 7" "460÷9=51, 1"
Replace-Text "541÷4=135, 1" "664÷6=110, 4"
Replace-Text "342÷7=48, 6" "852÷2=426, 0"

Replace-Text "944÷9=104, 8" "483÷7=69, 0"
Replace-Text "115÷3=38, 1" "440÷5=88, 0"
Replace-Text "289÷7=41, 2" "561÷2=280, 1"
Replace-Text "950÷8=118, 6" "437÷9=48, 5"
Replace-Text "151÷8=18, 7" "636÷4=159, 0"

Replace-Text "434÷4=108, 2" "145÷9=16, 1"
Replace-Text "769÷8=96, 1" "382÷3=127, 1"
Replace-Text "244÷4=61, 0" "578÷7=82, 4"
Replace-Text "930÷7=132, 6" "255÷3=85, 0"
Replace-Text "847÷6=141, 1" "915÷6=152, 3"

Replace-Text "152÷6=25, 2" "209÷3=69, 2"
Replace-Text "869÷4=217, 1" "304÷7=43, 3"
Replace-Text "138÷2=69, 0" "185÷5=37, 0"
Replace-Text "836÷4=209, 0" "667÷9=74, 1"
Replace-Text "995÷8=124, 3" "806÷2=403, 0"

Write-Output "Replacements complete"
